# Auto-generated edit script applying cell-level updates to cryptos sheet
# Updated symbol list with new coin prices/volumes (per commit message)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = "'323.04"
$ws.Range('E2').Value = "'-2.77%"
$ws.Range('D3').Value = "'43.29"
$ws.Range('E3').Value = "'-5.43%"
$ws.Range('D4').Value = "'5.276"
$ws.Range('E4').Value = "'-7.14%"
$ws.Range('E5').Value = "'-2.58%"
$ws.Range('B6').Value = 'GateToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D6').Value = "'4.383"
$ws.Range('E6').Value = "'-1.77%"
$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D7').Value = "'1.798"
$ws.Range('E7').Value = "'-12.00%"
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').Value = "'0.9618"
$ws.Range('E8').Value = "'-1.92%"
$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D9').Value = "'0.1121"
$ws.Range('E9').Value = "'-3.55%"
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = "'0.1864"
$ws.Range('E10').Value = "'-3.85%"
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = "'0.09419"
$ws.Range('E11').Value = "'-6.21%"
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').Value = "'0.04615"
$ws.Range('E12').Value = "'-0.27%"
$ws.Range('B13').Value = 'MCDex'
$ws.Range('C13').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D13').Value = "'7.494"
$ws.Range('E13').Value = "'-27.74%"
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = "'0.1064"
$ws.Range('E14').Value = "'0.44%"
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = "'0.001286"
$ws.Range('E15').Value = "'0.52%"
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = "'0.005828"
$ws.Range('E16').Value = "'-4.33%"
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = "'3.370"
$ws.Range('E17').Value = "'0.00%"
$ws.Range('D18').Value = "'2.517"
$ws.Range('E18').Value = "'-3.08%"
$ws.Range('D19').Value = "'0.3362"
$ws.Range('E19').Value = "'0.36%"
$ws.Range('E20').Value = "'-0.62%"
$ws.Range('D21').Value = "'0.2594"
$ws.Range('E21').Value = "'0.00%"
$ws.Range('D22').Value = "'0.04221"
$ws.Range('E22').Value = "'0.22%"
$ws.Range('D23').Value = "'0.001263"
$ws.Range('E23').Value = "'-3.46%"
$ws.Range('D24').Value = "'0.004312"
$ws.Range('E24').Value = "'-6.29%"
$ws.Range('D25').Value = "'0.0001311"
$ws.Range('E25').Value = "'2.34%"
$ws.Range('D26').Value = "'0.0003004"
$ws.Range('E26').Value = "'-19.71%"
$ws.Range('D38').Value = "'0.02641"
$ws.Range('E38').Value = "'-4.82%"
$ws.Range('D39').Value = "'0.05495"
$ws.Range('E39').Value = "'-5.41%"
$ws.Range('D40').Value = "'0.007854"
$ws.Range('E40').Value = "'1.59%"
$ws.Range('D41').Value = "'0.1400"
$ws.Range('E41').Value = "'-2.66%"
$ws.Range('D42').Value = "'0.006597"
$ws.Range('E42').Value = "'-8.34%"
$ws.Range('D43').Value = "'0.002131"
$ws.Range('E43').Value = "'7.93%"
$ws.Range('D44').Value = "'0.008696"
$ws.Range('E44').Value = "'6.41%"
$ws.Range('D45').Value = "'0.3305"
$ws.Range('D46').Value = "'0.00007026"
$ws.Range('E46').Value = "'-2.42%"
$ws.Range('D47').Value = "'0.00000000756"
$ws.Range('E47').Value = "'0.79%"
$ws.Range('D48').Value = "'0.003497"
$ws.Range('E48').Value = "'0.13%"
$ws.Range('D49').Value = "'0.003560"
$ws.Range('E49').Value = "'1.72%"
$ws.Range('D50').Value = "'0.00002118"
$ws.Range('E50').Value = "'0.79%"
$ws.Range('D51').Value = "'0.0002017"
$ws.Range('E51').Value = "'0.79%"
